# === Edit script: apply the HPC-ED.xlsx changes ===
$wb = $excel.ActiveWorkbook

# --- 1. Add the new hidden "DD" (data dictionary) sheet at the end ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$dd = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$dd.Name = "DD"

# --- 2. Populate the DD sheet with the lookup values used by the validations ---
$ddData = New-Object 'object[,]' 20,9
$ddData[0,0] = 'id'
$ddData[0,1] = 'visible_to'
$ddData[0,2] = 'Resource_URL_Type'
$ddData[0,3] = 'Language'
$ddData[0,4] = 'Cost'
$ddData[0,5] = 'Expertise_Level'
$ddData[0,6] = 'Learning_Outcome'
$ddData[0,7] = 'License'
$ddData[0,8] = 'Target_Group'
$ddData[1,0] = 'std'
$ddData[1,1] = 'public'
$ddData[1,2] = 'URL'
$ddData[1,3] = 'en'
$ddData[1,4] = 'no'
$ddData[1,5] = 'Beginner'
$ddData[1,6] = 'Remember'
$ddData[1,7] = 'Creative Commons'
$ddData[1,8] = 'Researchers'
$ddData[2,2] = 'ARK'
$ddData[2,4] = 'yes'
$ddData[2,5] = 'Intermediate'
$ddData[2,6] = 'Understand'
$ddData[2,7] = 'MIT'
$ddData[2,8] = 'Research groups'
$ddData[3,2] = 'arXiv'
$ddData[3,4] = 'maybe'
$ddData[3,5] = 'Advanced'
$ddData[3,6] = 'Apply'
$ddData[3,7] = 'BSD'
$ddData[3,8] = 'Research communities'
$ddData[4,2] = 'bibcode'
$ddData[4,5] = 'All'
$ddData[4,6] = 'Analyze'
$ddData[4,7] = 'GPL'
$ddData[4,8] = 'Research projects'
$ddData[5,2] = 'DOI'
$ddData[5,6] = 'Evaluate'
$ddData[5,7] = 'GNU'
$ddData[5,8] = 'Research networks'
$ddData[6,2] = 'EAN13'
$ddData[6,6] = 'Create'
$ddData[6,8] = 'Research managers'
$ddData[7,2] = 'EISSN'
$ddData[7,8] = 'Research organizations'
$ddData[8,2] = 'Handle'
$ddData[8,8] = 'Student'
$ddData[9,2] = 'IGSN'
$ddData[9,8] = 'Innovators'
$ddData[10,2] = 'ISBN'
$ddData[10,8] = 'Providers'
$ddData[11,2] = 'ISSN'
$ddData[11,8] = 'Funders'
$ddData[12,2] = 'ISTC'
$ddData[12,8] = 'Research Infrastructure Managers'
$ddData[13,2] = 'LISSN'
$ddData[13,8] = 'Resource Managers'
$ddData[14,2] = 'LSID'
$ddData[14,8] = 'Publishers'
$ddData[15,2] = 'PMID'
$ddData[15,8] = 'Other'
$ddData[16,2] = 'PURL'
$ddData[17,2] = 'UPC'
$ddData[18,2] = 'URN'
$ddData[19,2] = 'w3id'
$dd.Range("A1:I20").Value = $ddData

# Row height / wrap tweaks that mirror the authored sheet
$dd.Rows.Item(8).RowHeight = 16
$dd.Rows.Item(9).RowHeight = 17
$dd.Range("C9").WrapText = $true

# Column widths
$dd.Columns.Item(1).ColumnWidth = 4.83203125
$dd.Columns.Item(2).ColumnWidth = 8.6640625
$dd.Columns.Item(3).ColumnWidth = 18.1640625
$dd.Columns.Item(6).ColumnWidth = 17.5
$dd.Range("G1:H1").ColumnWidth = 18.6640625
$dd.Columns.Item(9).ColumnWidth = 19.5

$dd.Range("E13").Select()

# --- 3. Hide the DD sheet (it is only a lookup sheet for data validation) ---
$dd.Visible = $false

# --- 4. Switch to "HPC-ED Share" and add the header cell comments ---
$ws = $wb.Worksheets.Item("HPC-ED Share")
$ws.Activate()
$ws.Range("A1").AddComment('A unique identifier for your material. It should begin with your Provider_ID, a colon, then an identifier for your material. Please see "Required Metadata" for full information. 
') | Out-Null
$ws.Range("B1").AddComment('id is always std') | Out-Null
$ws.Range("C1").AddComment('visible_to is always public') | Out-Null
$ws.Range("D1").AddComment('The title of your material, e.g. Intro to Slurm') | Out-Null
$ws.Range("E1").AddComment('A link to your material.') | Out-Null
$ws.Range("F1").AddComment('The tpye of URL given. Most common choice is URL.') | Out-Null
$ws.Range("G1").AddComment('Language the material is written in. For English, it is en.') | Out-Null
$ws.Range("H1").AddComment('Whether the material charges for access. Typically, no for free materials.') | Out-Null
$ws.Range("I1").AddComment('Your assigned Provider_ID. It is unique to your organization.') | Out-Null
$ws.Range("J1").AddComment('The abstract of your material.') | Out-Null
$ws.Range("K1").AddComment('The author(s) of the material. Use commas to separate for multiple auithors.') | Out-Null
$ws.Range("L1").AddComment('The level of skill needed to complete the material.') | Out-Null
$ws.Range("M1").AddComment('The keywords that describe the material. What is the material about?') | Out-Null
$ws.Range("N1").AddComment('The intended outcome of the material. What is your material hoping to achieve?') | Out-Null
$ws.Range("O1").AddComment('The type of material. Is it a recorded lesson? A lesson plan? Or a textbook?') | Out-Null
$ws.Range("P1").AddComment('The license of the material.') | Out-Null
$ws.Range("Q1").AddComment('The target audience of the material. Who is it intended for?') | Out-Null
$ws.Range("R1").AddComment('The release date of the material, in ISO 8601 format. For example, 2023-11-05T08:15:30-05:00 for November 5, 2023, 8:15:30 am, US Eastern Standard Time. ') | Out-Null
$ws.Range("S1").AddComment('The start date of the material, in ISO 8601 format. For example, 2023-11-05T08:15:30-05:00 for November 5, 2023, 8:15:30 am, US Eastern Standard Time. 
') | Out-Null
$ws.Range("T1").AddComment('Estimated time to complete this material in minutes. ') | Out-Null
$ws.Range("U1").AddComment('From 0.0 to 5.0, the rating of the material.') | Out-Null

# --- 5. Data validation rules on "HPC-ED Share" ---
$val = $ws.Range("T1:T1048576").Validation
$val.Add(1, 2, 7, "0")
$val.ErrorTitle = "invliad Duration"
$val.ShowInput = $false
$val.ShowError = $true

$val = $ws.Range("U1:U1048576").Validation
$val.Add(2, 2, 1, "0", "5")
$val.ErrorTitle = "invalid Rating"
$val.ShowInput = $false
$val.ShowError = $true

$val = $ws.Range("B1:B1048576").Validation
$val.Add(3, 2, 1, '=DD!\$A\$2')
$val.ErrorTitle = "id must be std"
$val.ShowInput = $false
$val.ShowError = $true

$val = $ws.Range("C1:C1048576").Validation
$val.Add(3, 2, 1, '=DD!\$B\$2')
$val.ShowInput = $false
$val.ShowError = $true

$val = $ws.Range("G1:G1048576").Validation
$val.Add(3, 1, 1, '=DD!\$D\$2')
$val.ShowInput = $false
$val.ShowError = $false

$val = $ws.Range("H1:H1048576").Validation
$val.Add(3, 1, 1, '=DD!\$E\$2:\$E\$4')
$val.ShowInput = $false
$val.ShowError = $false

$val = $ws.Range("L1:L1048576").Validation
$val.Add(3, 1, 1, '=DD!\$F\$2:\$F\$5')
$val.ShowInput = $false
$val.ShowError = $false

$val = $ws.Range("N1:N1048576").Validation
$val.Add(3, 1, 1, '=DD!\$G\$2:\$G\$7')
$val.ShowInput = $false
$val.ShowError = $false

$val = $ws.Range("P1:P1048576").Validation
$val.Add(3, 1, 1, '=DD!\$H\$2:\$H\$6')
$val.ShowInput = $false
$val.ShowError = $false

$val = $ws.Range("Q1:Q1048576").Validation
$val.Add(3, 1, 1, '=DD!\$I\$2:\$I\$16')
$val.ShowInput = $false
$val.ShowError = $false

$val = $ws.Range("F1:F1048576").Validation
$val.Add(3, 1, 1, '=DD!\$C\$2:\$C\$20')
$val.ShowInput = $false
$val.ShowError = $false

# --- 6. Restore the original selection on this sheet, then make it the active tab ---
$ws.Range("U6").Select()

# --- 7. Update the selection remembered on the "Example" sheet ---
$example = $wb.Worksheets.Item("Example")
$example.Range("D24").Select()

# --- 8. Make sure "HPC-ED Share" ends up as the active sheet/tab ---
$ws.Activate()
$ws.Range("U6").Select()
